$d = $word.ActiveDocument

# --- Title ---
$d.Paragraphs(1).Range.Find.Execute(
    "The Quantum Realm Unlocked: Exploring the Enigmatic World of Subatomic Particles",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Exploring the Marvels of Mathematics: A Journey Through Numbers and Patterns", 2)

# --- Author name line: "Dr. Aria Walker" -> "Ms. Emily Watson" ---
$d.Paragraphs(2).Range.Find.Execute(
    "Dr", $true, $false, $false, $false, $false, $true, 1, $false, "Ms", 2)
$d.Paragraphs(2).Range.Find.Execute(
    " Aria Walker", $true, $false, $false, $false, $false, $true, 1, $false, " Emily Watson", 2)

# --- Contact line: "walker.aria@quantumstudies.org" -> "at" ---
$d.Paragraphs(3).Range.Find.Execute(
    "walker.aria@quantumstudies.org", $true, $false, $false, $false, $false, $true, 1, $false, "at", 2)

# --- Body paragraph (long paragraph after the blank line) ---
$body = $d.Paragraphs(5).Range
$body.Find.Execute(
    "In the realm of science, few domains are more captivating and enigmatic than the quantum world",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Mathematics, a subject often shrouded in complexity, is an art form that unveils the hidden beauty of numbers", 2)

$body = $d.Paragraphs(5).Range
$body.Find.Execute(
    " This intricate and mysterious universe of subatomic entities, characterized by phenomena such as superposition, entanglement, and wave-particle duality, belies our classical understanding of reality",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " It is a universal language that has captivated countless minds throughout history, from ancient civilizations to modern-day scientists", 2)

$body = $d.Paragraphs(5).Range
$body.Find.Execute(
    " In seeking to unveil the secrets of the quantum realm, scientists embark on a thrilling journey to unravel the fundamental fabric of the universe, unlocking new paradigms of knowledge and innovation. As we delve deeper into the quantum realm, we encounter particles that behave in ways akin to waves, blurring the lines between distinct states. This enigmatic phenomenon, known as wave-particle duality, challenges our classical notions and compels us to reconceptualize the very nature of matter and energy",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " This journey into the world of mathematics takes us through a mesmerizing tapestry of patterns, shapes, and equations, unravelling the enigmatic mysteries that lie beneath the surface", 2)

$body = $d.Paragraphs(5).Range
$body.Find.Execute(
    "Moreover, quantum entanglement serves as a perplexing paradox that defies explanation within classical physics",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Mathematics is an intricate dance of numbers, a symphony of symbols that plays out on the stage of our minds", 2)

$body = $d.Paragraphs(5).Range
$body.Find.Execute(
    " In this intricate dance of subatomic particles, regardless of the distance separating them, the actions of one influence the state of the others",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " It is a language of logic and reason, allowing us to understand the complexities of the world around us", 2)

$body = $d.Paragraphs(5).Range
$body.Find.Execute(
    " This interconnectedness transcends the constraints of space and time, raising profound questions about locality, causality, and the interconnectedness of all things",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Its elegance lies in its simplicity, yet its depth reveals an inexhaustible well of knowledge, intriguing puzzles, and awe-inspiring discoveries", 2)

$body = $d.Paragraphs(5).Range
$body.Find.Execute(
    " Furthermore, the concept of superposition places particles in an ambiguous state, simultaneously existing in multiple states or locations until observed, a phenomenon that further challenges our conventional understanding of the physical world",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "^l^lAs we delve deeper into the realm of mathematics, we discover the interconnectedness of concepts, the harmony between theory and application", 2)

$body = $d.Paragraphs(5).Range
$body.Find.Execute(
    " Unraveling these enigmas promises transformative insights with wide-ranging implications, ranging from the development of quantum computers, surpassing classical computational limits, to the establishment of secure communication protocols leveraging quantum information theory",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " From the Pythagorean theorem to the calculus of infinitesimals, mathematics provides a framework for understanding phenomena across diverse fields, from engineering to finance to the boundless mysteries of the cosmos.^l^lBody:^l^lMathematics, in its essence, is a journey of exploration and discovery. It invites us to question, to seek answers, to push the boundaries of human knowledge. Through the lens of mathematics, we unravel the secrets of nature, unlocking the mysteries of the universe. Equations become tools that decipher the intricate dance of particles, the symphony of celestial bodies, and the intricate patterns of genetic codes.^l^lThis pursuit of knowledge is a collaborative effort, a collective exploration of the human race. We stand on the shoulders of giants, building upon the legacy of those who came before us. Through textbooks and classrooms, we inherit the wisdom of countless mathematicians, each contributing a piece to the ever-expanding mosaic of knowledge.^l^lMathematics is not merely a subject to be studied; it is an art to be appreciated, a skill to be mastered. It is a mindset, a way of thinking that permeates every aspect of our lives. The analytical thinking honed through mathematical practice extends far beyond the classroom, providing a valuable tool for decision-making, problem-solving, and navigating the complexities of the world", 2)

Write-Output "----"
Write-Output $d.Paragraphs(5).Range.Text

# --- Summary paragraph (paragraph 7) ---
$summary = $d.Paragraphs(7).Range
$summary.Find.Execute(
    "In conclusion, the quest to comprehend the enigmatic quantum realm undeniably represents one of science's most profound and awe-inspiring endeavors",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In conclusion, the world of mathematics is a vast and awe-inspiring landscape, an ever-evolving testament to human ingenuity", 2)

$summary = $d.Paragraphs(7).Range
$summary.Find.Execute(
    " As physicists unravel the secrets of subatomic particles, they unlock gateways to understanding the very essence of matter and energy",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Its beauty lies in its simplicity and elegance, yet its depth reveals an inexhaustible wellspring of knowledge", 2)

$summary = $d.Paragraphs(7).Range
$summary.Find.Execute(
    " Driven by an insatiable thirst for knowledge, scientists continue their voyage into the quantum realm, illuminated by the beacons of superposition, entanglement, and wave-particle duality",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " As we embark on this mathematical journey, we discover the interconnectedness of concepts, the harmony between theory and application, and the boundless potential of human understanding", 2)

$summary = $d.Paragraphs(7).Range
$summary.Find.Execute(
    " Along this extraordinary journey, they not only decipher the enigmas of the unseen but also pave the way for transformative technologies that promise to reshape our future",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Mathematics is more than just a subject; it is a mindset, a skill, and an art form that shapes our perception of the world and empowers us to make sense of its complexities", 2)

$summary = $d.Paragraphs(7).Range
$summary.Find.Execute(
    " The quantum world unveils a realm where the laws of classical physics break down, giving way to a symphony of paradoxical phenomena and infinite possibilities",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " As we continue to explore this realm of numbers and patterns, we unlock the secrets of the universe and embark on a lifelong journey of discovery", 2)

Write-Output "----"
Write-Output $d.Paragraphs(7).Range.Text

# --- Add a trailing empty paragraph at the end of the document body ---
$summary = $d.Paragraphs(7).Range
$summary.Find.Execute(
    "embark on a lifelong journey of discovery.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "embark on a lifelong journey of discovery.^p", 2)

Write-Output "----"
Write-Output $d.Paragraphs.Count
